# Fix dataset issues in "commenti_per_topic" sheet:
#  - row 3 topic name had a typo introduced (CRONACA -> CROANCA)
#  - a missing data row (POLITICA / Instagram) was inserted, shifting the
#    social values for several topic/social combinations down by one row
#  - negativo/positivo counts were corrected across the board

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# topic, social, negativo, positivo  for rows 3..12
$data = @(
    @("CROANCA",      "Facebook",   34,   6),
    @("CRONACA",      "Facebook",   806,  346),
    @("CRONACA",      "Instagram",  958,  179),
    @("CRONACA",      "YouTube",    906,  275),
    @("CRONACA NERA", "Facebook",   1064, 136),
    @("CRONACA NERA", "Instagram",  1044, 156),
    @("CRONACA NERA", "YouTube",    1023, 177),
    @("POLITICA",     "Facebook",   874,  325),
    @("POLITICA",     "Instagram",  967,  226),
    @("POLITICA",     "YouTube",    895,  300)
)

$row = 3
foreach ($entry in $data) {
    if ($row -eq 12) {
        # new row appended at the bottom; materialize the (empty) column A
        # cell the same way the rest of the table rows have it, without
        # introducing any new style
        $ws.Cells.Item($row, 1).Style = $ws.Cells.Item(3, 1).Style
    }
    $ws.Cells.Item($row, 2).Value = $entry[0]
    $ws.Cells.Item($row, 3).Value = $entry[1]
    $ws.Cells.Item($row, 4).Value = $entry[2]
    $ws.Cells.Item($row, 5).Value = $entry[3]
    $row = $row + 1
}
